$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("CRbQ")

# Rename the existing "hydrogen" power-plant-type row to "hydrogen combustion turbine"
$ws.Cells.Item(24, 1).Value = "hydrogen combustion turbine"

# Give it bold / vertically centered styling
$ws.Cells.Item(24, 1).Font.Bold = $true
$ws.Cells.Item(24, 1).Font.Color = 0
$ws.Cells.Item(24, 1).VerticalAlignment = -4108

# Add a new row for "hydrogen combined cycle" as a new power plant type, using
# the same formatting as the row above (copy format to avoid creating a stray
# unused cell style)
$ws.Cells.Item(25, 1).Value = "hydrogen combined cycle"
$ws.Cells.Item(24, 1).Copy() | Out-Null
$ws.Cells.Item(25, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range($ws.Cells.Item(25, 2), $ws.Cells.Item(25, 32)).Value = 0

# Update the on-screen selection to the new row's data range, then restore
# the "About" sheet as the active tab (matches saved workbook state)
$ws.Range("B25:AF25").Select() | Out-Null
$wsAbout.Activate() | Out-Null
